$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (AB1, AC1) - new question columns, matching header style (bold, centered) of existing header cells
$ws.Cells.Item(1, 28).Value = "Q_Responsibility"
$ws.Cells.Item(1, 28).Font.Bold = $true
$ws.Cells.Item(1, 28).HorizontalAlignment = $ws.Cells.Item(1, 27).HorizontalAlignment
$ws.Cells.Item(1, 29).Value = "Q_Responsibilitycode"
$ws.Cells.Item(1, 29).Font.Bold = $true
$ws.Cells.Item(1, 29).HorizontalAlignment = $ws.Cells.Item(1, 27).HorizontalAlignment

# Data rows: AB = Q_Responsibility (text), AC = Q_Responsibilitycode (numeric-looking code stored as text)
$ws.Cells.Item(2, 28).Value = "Public authorities are responsible and citizens somewhat responsible for flood protection"
$ws.Cells.Item(2, 29).NumberFormat = "@"
$ws.Cells.Item(2, 29).Value = "2"
$ws.Cells.Item(2, 29).ClearFormats()
$ws.Cells.Item(3, 28).Value = "Public authorities and citizens are equally responsible for flood protection"
$ws.Cells.Item(3, 29).NumberFormat = "@"
$ws.Cells.Item(3, 29).Value = "3"
$ws.Cells.Item(3, 29).ClearFormats()
$ws.Cells.Item(4, 28).Value = "Public authorities are completely responsible for flood protection"
$ws.Cells.Item(4, 29).NumberFormat = "@"
$ws.Cells.Item(4, 29).Value = "1"
$ws.Cells.Item(4, 29).ClearFormats()
$ws.Cells.Item(5, 28).Value = "Public authorities are responsible and citizens somewhat responsible for flood protection"
$ws.Cells.Item(5, 29).NumberFormat = "@"
$ws.Cells.Item(5, 29).Value = "2"
$ws.Cells.Item(5, 29).ClearFormats()
$ws.Cells.Item(6, 28).Value = "Public authorities are completely responsible for flood protection"
$ws.Cells.Item(6, 29).NumberFormat = "@"
$ws.Cells.Item(6, 29).Value = "1"
$ws.Cells.Item(6, 29).ClearFormats()
$ws.Cells.Item(7, 28).Value = "Public authorities are responsible and citizens somewhat responsible for flood protection"
$ws.Cells.Item(7, 29).NumberFormat = "@"
$ws.Cells.Item(7, 29).Value = "2"
$ws.Cells.Item(7, 29).ClearFormats()
$ws.Cells.Item(8, 28).Value = "Public authorities are completely responsible for flood protection"
$ws.Cells.Item(8, 29).NumberFormat = "@"
$ws.Cells.Item(8, 29).Value = "1"
$ws.Cells.Item(8, 29).ClearFormats()
$ws.Cells.Item(9, 28).Value = "Public authorities are responsible and citizens somewhat responsible for flood protection"
$ws.Cells.Item(9, 29).NumberFormat = "@"
$ws.Cells.Item(9, 29).Value = "2"
$ws.Cells.Item(9, 29).ClearFormats()
$ws.Cells.Item(10, 28).Value = "Public authorities are responsible and citizens somewhat responsible for flood protection"
$ws.Cells.Item(10, 29).NumberFormat = "@"
$ws.Cells.Item(10, 29).Value = "2"
$ws.Cells.Item(10, 29).ClearFormats()
$ws.Cells.Item(11, 28).Value = "Public authorities are completely responsible for flood protection"
$ws.Cells.Item(11, 29).NumberFormat = "@"
$ws.Cells.Item(11, 29).Value = "1"
$ws.Cells.Item(11, 29).ClearFormats()
$ws.Cells.Item(12, 28).Value = "Public authorities are responsible and citizens somewhat responsible for flood protection"
$ws.Cells.Item(12, 29).NumberFormat = "@"
$ws.Cells.Item(12, 29).Value = "2"
$ws.Cells.Item(12, 29).ClearFormats()
$ws.Cells.Item(13, 28).Value = "Public authorities are completely responsible for flood protection"
$ws.Cells.Item(13, 29).NumberFormat = "@"
$ws.Cells.Item(13, 29).Value = "1"
$ws.Cells.Item(13, 29).ClearFormats()
$ws.Cells.Item(14, 28).Value = "Public authorities and citizens are equally responsible for flood protection"
$ws.Cells.Item(14, 29).NumberFormat = "@"
$ws.Cells.Item(14, 29).Value = "3"
$ws.Cells.Item(14, 29).ClearFormats()
$ws.Cells.Item(15, 28).Value = "Public authorities are completely responsible for flood protection"
$ws.Cells.Item(15, 29).NumberFormat = "@"
$ws.Cells.Item(15, 29).Value = "1"
$ws.Cells.Item(15, 29).ClearFormats()
$ws.Cells.Item(16, 28).Value = "Public authorities are responsible and citizens somewhat responsible for flood protection"
$ws.Cells.Item(16, 29).NumberFormat = "@"
$ws.Cells.Item(16, 29).Value = "2"
$ws.Cells.Item(16, 29).ClearFormats()
$ws.Cells.Item(17, 28).Value = "Public authorities and citizens are equally responsible for flood protection"
$ws.Cells.Item(17, 29).NumberFormat = "@"
$ws.Cells.Item(17, 29).Value = "3"
$ws.Cells.Item(17, 29).ClearFormats()
$ws.Cells.Item(18, 28).Value = "Public authorities are responsible and citizens somewhat responsible for flood protection"
$ws.Cells.Item(18, 29).NumberFormat = "@"
$ws.Cells.Item(18, 29).Value = "2"
$ws.Cells.Item(18, 29).ClearFormats()
$ws.Cells.Item(19, 28).Value = "Public authorities are responsible and citizens somewhat responsible for flood protection"
$ws.Cells.Item(19, 29).NumberFormat = "@"
$ws.Cells.Item(19, 29).Value = "2"
$ws.Cells.Item(19, 29).ClearFormats()
$ws.Cells.Item(20, 28).Value = "Public authorities are completely responsible for flood protection"
$ws.Cells.Item(20, 29).NumberFormat = "@"
$ws.Cells.Item(20, 29).Value = "1"
$ws.Cells.Item(20, 29).ClearFormats()
$ws.Cells.Item(21, 28).Value = "Public authorities are completely responsible for flood protection"
$ws.Cells.Item(21, 29).NumberFormat = "@"
$ws.Cells.Item(21, 29).Value = "1"
$ws.Cells.Item(21, 29).ClearFormats()
$ws.Cells.Item(22, 28).Value = "Public authorities and citizens are equally responsible for flood protection"
$ws.Cells.Item(22, 29).NumberFormat = "@"
$ws.Cells.Item(22, 29).Value = "3"
$ws.Cells.Item(22, 29).ClearFormats()
$ws.Cells.Item(23, 28).Value = "Public authorities and citizens are equally responsible for flood protection"
$ws.Cells.Item(23, 29).NumberFormat = "@"
$ws.Cells.Item(23, 29).Value = "3"
$ws.Cells.Item(23, 29).ClearFormats()
$ws.Cells.Item(24, 28).Value = "Public authorities are responsible and citizens somewhat responsible for flood protection"
$ws.Cells.Item(24, 29).NumberFormat = "@"
$ws.Cells.Item(24, 29).Value = "2"
$ws.Cells.Item(24, 29).ClearFormats()
$ws.Cells.Item(25, 28).Value = "Public authorities are responsible and citizens somewhat responsible for flood protection"
$ws.Cells.Item(25, 29).NumberFormat = "@"
$ws.Cells.Item(25, 29).Value = "2"
$ws.Cells.Item(25, 29).ClearFormats()
$ws.Cells.Item(26, 28).Value = "Public authorities are responsible and citizens somewhat responsible for flood protection"
$ws.Cells.Item(26, 29).NumberFormat = "@"
$ws.Cells.Item(26, 29).Value = "2"
$ws.Cells.Item(26, 29).ClearFormats()
$ws.Cells.Item(27, 28).Value = "Public authorities are completely responsible for flood protection"
$ws.Cells.Item(27, 29).NumberFormat = "@"
$ws.Cells.Item(27, 29).Value = "1"
$ws.Cells.Item(27, 29).ClearFormats()
$ws.Cells.Item(29, 28).Value = "Public authorities are responsible and citizens somewhat responsible for flood protection"
$ws.Cells.Item(29, 29).NumberFormat = "@"
$ws.Cells.Item(29, 29).Value = "2"
$ws.Cells.Item(29, 29).ClearFormats()
$ws.Cells.Item(30, 28).Value = "Public authorities are completely responsible for flood protection"
$ws.Cells.Item(30, 29).NumberFormat = "@"
$ws.Cells.Item(30, 29).Value = "1"
$ws.Cells.Item(30, 29).ClearFormats()
$ws.Cells.Item(31, 28).Value = "Public authorities are completely responsible for flood protection"
$ws.Cells.Item(31, 29).NumberFormat = "@"
$ws.Cells.Item(31, 29).Value = "1"
$ws.Cells.Item(31, 29).ClearFormats()
$ws.Cells.Item(32, 28).Value = "Public authorities are responsible and citizens somewhat responsible for flood protection"
$ws.Cells.Item(32, 29).NumberFormat = "@"
$ws.Cells.Item(32, 29).Value = "2"
$ws.Cells.Item(32, 29).ClearFormats()
$ws.Cells.Item(33, 28).Value = "Public authorities are responsible and citizens somewhat responsible for flood protection"
$ws.Cells.Item(33, 29).NumberFormat = "@"
$ws.Cells.Item(33, 29).Value = "2"
$ws.Cells.Item(33, 29).ClearFormats()
$ws.Cells.Item(34, 28).Value = "Public authorities are responsible and citizens somewhat responsible for flood protection"
$ws.Cells.Item(34, 29).NumberFormat = "@"
$ws.Cells.Item(34, 29).Value = "2"
$ws.Cells.Item(34, 29).ClearFormats()
$ws.Cells.Item(35, 28).Value = "Public authorities are responsible and citizens somewhat responsible for flood protection"
$ws.Cells.Item(35, 29).NumberFormat = "@"
$ws.Cells.Item(35, 29).Value = "2"
$ws.Cells.Item(35, 29).ClearFormats()
$ws.Cells.Item(36, 28).Value = "Public authorities are responsible and citizens somewhat responsible for flood protection"
$ws.Cells.Item(36, 29).NumberFormat = "@"
$ws.Cells.Item(36, 29).Value = "2"
$ws.Cells.Item(36, 29).ClearFormats()
$ws.Cells.Item(37, 28).Value = "Public authorities are completely responsible for flood protection"
$ws.Cells.Item(37, 29).NumberFormat = "@"
$ws.Cells.Item(37, 29).Value = "1"
$ws.Cells.Item(37, 29).ClearFormats()
$ws.Cells.Item(38, 28).Value = "Public authorities are responsible and citizens somewhat responsible for flood protection"
$ws.Cells.Item(38, 29).NumberFormat = "@"
$ws.Cells.Item(38, 29).Value = "2"
$ws.Cells.Item(38, 29).ClearFormats()
$ws.Cells.Item(39, 28).Value = "Public authorities are responsible and citizens somewhat responsible for flood protection"
$ws.Cells.Item(39, 29).NumberFormat = "@"
$ws.Cells.Item(39, 29).Value = "2"
$ws.Cells.Item(39, 29).ClearFormats()
$ws.Cells.Item(40, 28).Value = "Public authorities and citizens are equally responsible for flood protection"
$ws.Cells.Item(40, 29).NumberFormat = "@"
$ws.Cells.Item(40, 29).Value = "3"
$ws.Cells.Item(40, 29).ClearFormats()
$ws.Cells.Item(41, 28).Value = "Public authorities are completely responsible for flood protection"
$ws.Cells.Item(41, 29).NumberFormat = "@"
$ws.Cells.Item(41, 29).Value = "1"
$ws.Cells.Item(41, 29).ClearFormats()
$ws.Cells.Item(42, 28).Value = "Public authorities are responsible and citizens somewhat responsible for flood protection"
$ws.Cells.Item(42, 29).NumberFormat = "@"
$ws.Cells.Item(42, 29).Value = "2"
$ws.Cells.Item(42, 29).ClearFormats()
$ws.Cells.Item(43, 28).Value = "Public authorities are completely responsible for flood protection"
$ws.Cells.Item(43, 29).NumberFormat = "@"
$ws.Cells.Item(43, 29).Value = "1"
$ws.Cells.Item(43, 29).ClearFormats()
$ws.Cells.Item(44, 28).Value = "Public authorities are responsible and citizens somewhat responsible for flood protection"
$ws.Cells.Item(44, 29).NumberFormat = "@"
$ws.Cells.Item(44, 29).Value = "2"
$ws.Cells.Item(44, 29).ClearFormats()
$ws.Cells.Item(45, 28).Value = "Public authorities and citizens are equally responsible for flood protection"
$ws.Cells.Item(45, 29).NumberFormat = "@"
$ws.Cells.Item(45, 29).Value = "3"
$ws.Cells.Item(45, 29).ClearFormats()
$ws.Cells.Item(46, 28).Value = "Public authorities are responsible and citizens somewhat responsible for flood protection"
$ws.Cells.Item(46, 29).NumberFormat = "@"
$ws.Cells.Item(46, 29).Value = "2"
$ws.Cells.Item(46, 29).ClearFormats()
$ws.Cells.Item(47, 28).Value = "Public authorities and citizens are equally responsible for flood protection"
$ws.Cells.Item(47, 29).NumberFormat = "@"
$ws.Cells.Item(47, 29).Value = "3"
$ws.Cells.Item(47, 29).ClearFormats()
$ws.Cells.Item(48, 28).Value = "Public authorities are responsible and citizens somewhat responsible for flood protection"
$ws.Cells.Item(48, 29).NumberFormat = "@"
$ws.Cells.Item(48, 29).Value = "2"
$ws.Cells.Item(48, 29).ClearFormats()
$ws.Cells.Item(49, 28).Value = "Public authorities are responsible and citizens somewhat responsible for flood protection"
$ws.Cells.Item(49, 29).NumberFormat = "@"
$ws.Cells.Item(49, 29).Value = "2"
$ws.Cells.Item(49, 29).ClearFormats()
$ws.Cells.Item(50, 28).Value = "Public authorities and citizens are equally responsible for flood protection"
$ws.Cells.Item(50, 29).NumberFormat = "@"
$ws.Cells.Item(50, 29).Value = "3"
$ws.Cells.Item(50, 29).ClearFormats()
$ws.Cells.Item(51, 28).Value = "Public authorities are responsible and citizens somewhat responsible for flood protection"
$ws.Cells.Item(51, 29).NumberFormat = "@"
$ws.Cells.Item(51, 29).Value = "2"
$ws.Cells.Item(51, 29).ClearFormats()
$ws.Cells.Item(52, 28).Value = "Public authorities and citizens are equally responsible for flood protection"
$ws.Cells.Item(52, 29).NumberFormat = "@"
$ws.Cells.Item(52, 29).Value = "3"
$ws.Cells.Item(52, 29).ClearFormats()
$ws.Cells.Item(53, 28).Value = "Public authorities are responsible and citizens somewhat responsible for flood protection"
$ws.Cells.Item(53, 29).NumberFormat = "@"
$ws.Cells.Item(53, 29).Value = "2"
$ws.Cells.Item(53, 29).ClearFormats()
